# Auto-generated: apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.397.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.45%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.905.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'526.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.94%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.42%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.28%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.915.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -4.73%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -7.56%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.72%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.412.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.26%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'60.469.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.57%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'22.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.895.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.74%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.0000140"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.88%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.96%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -3.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'360.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -8.66%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.74%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.59%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'63.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.43%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.012.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.448"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.53%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.175"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.91%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.08%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -9.36%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0851"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -11.50%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -4.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'19.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.84%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'151.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.89%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -7.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -7.92%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -8.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -7.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'37.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.96%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.336.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Mantle"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.646"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.41%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Filecoin"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'3.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'20.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.53%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0569"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.28%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.06%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'4.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'10.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -5.67%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.01%  "
$ws.Range("E51").Style = "Normal"
